$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: clear the stray C2 / E2 values (bug produced extra y_0_forecast /
# y_1_forecast entries that shouldn't exist for this row) ---
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# --- Row 3: clear stray C3, and correct the recomputed E3 forecast ---
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 4.422525088127305

# --- Row 4: corrected forecast values after bug fix ---
$ws.Range("C4").Value = -14.45332333832744
$ws.Range("E4").Value = -2.928447329610051

# --- Row 5 ---
$ws.Range("C5").Value = 8.600536527919612

# --- Row 7 ---
$ws.Range("C7").Value = 4.639893381363192

# --- Row 8 ---
$ws.Range("E8").Value = 2.429116709932599

# --- Row 9 ---
$ws.Range("E9").Value = 3.941300050092877

# --- Row 12 ---
$ws.Range("C12").Value = 4.695933104194361

# --- Row 15 ---
$ws.Range("E15").Value = -1.352810423674367

# --- Row 16 ---
$ws.Range("E16").Value = -0.3934198590721305

# --- Row 17 ---
$ws.Range("C17").Value = 5.120680133083622

# --- Row 18 ---
$ws.Range("C18").Value = -0.5532735011319123

# --- Row 19 ---
$ws.Range("E19").Value = -1.28528149926006
